$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-14 Saturday" "2025-06-15 Sunday"

Replace-Text "835×5=" "165×8="
Replace-Text "614×4=" "331×6="
Replace-Text "931×3=" "823×5="
Replace-Text "556×4=" "913×6="
Replace-Text "350×3=" "971×8="
Replace-Text "966×3=" "834×3="
Replace-Text "396×5=" "603×8="
Replace-Text "136×4=" "583×9="
Replace-Text "754×3=" "622×6="
Replace-Text "217×9=" "649×4="
Replace-Text "257×5=" "586×2="
Replace-Text "233×7=" "789×5="
Replace-Text "265×6=" "477×9="
Replace-Text "690×7=" "453×8="
Replace-Text "384×6=" "677×8="
Replace-Text "160×9=" "236×3="
Replace-Text "733×2=" "657×8="
Replace-Text "972×2=" "678×3="
Replace-Text "330×6=" "905×4="
Replace-Text "943×4=" "816×2="
Replace-Text "866×5=" "468×7="
Replace-Text "152×2=" "708×9="
Replace-Text "153×6=" "587×9="
Replace-Text "883×4=" "377×3="
Replace-Text "335×7=" "307×9="
